# Rename the embedded picture drawing objects in the footers/header.
#
#   * Primary footer   (PearsonLogo, docPr id="2")  image2.png -> image1.png
#   * First-page footer(PearsonLogo, docPr id="3")  image2.png -> image1.png
#   * First-page header(BTec_Logo-Orange, docPr id="1") image1.jpg -> image2.jpg
#
# WdHeaderFooterIndex constants aren't predefined in this host, so the
# literal values are used: 1 = wdHeaderFooterPrimary, 2 = wdHeaderFooterFirstPage.

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# --- Primary footer: PearsonLogo (image2.png -> image1.png) ---
$pearsonPrimary = $section.Footers.Item(1).Range.InlineShapes.Item(1)
$pearsonPrimary.Range.InlineShapes.Item(1).Name = "image1.png"

# --- First-page footer: PearsonLogo (image2.png -> image1.png) ---
$pearsonFirst = $section.Footers.Item(2).Range.InlineShapes.Item(1)
$pearsonFirst.Range.InlineShapes.Item(1).Name = "image1.png"

# --- First-page header: BTec_Logo-Orange (image1.jpg -> image2.jpg) ---
$btecFirst = $section.Headers.Item(2).Range.InlineShapes.Item(1)
$btecFirst.Range.InlineShapes.Item(1).Name = "image2.jpg"

Write-Host "Renamed picture drawing objects."
